# Rotate the three mushroom-find records currently on rows 25, 27, 28:
#   new row25 <- old row27 data (plus a location comment + biotope info)
#   new row27 <- old row28 data (plus an updated location comment)
#   new row28 <- old row25 data (with its location comment / biotope cleared)
# Row 26 is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: becomes the old row-27 record, with new/added fields ---
$ws.Cells.Item(25, 1).Value  = 112306179            # A25 Id
$ws.Cells.Item(25, 2).Value  = 89114                # B25 Taxonsorteringsordning
$ws.Cells.Item(25, 4).Value  = "NT"                 # D25 Rodlistade
$ws.Cells.Item(25, 5).Value  = 5754                 # E25 TaxonId
$ws.Cells.Item(25, 6).Value  = "Gultoppig fingersvamp"   # F25 Artnamn
$ws.Cells.Item(25, 7).Value  = "Ramaria testaceoflava"   # G25 Vetenskapligt namn
$ws.Cells.Item(25, 8).Value  = "(Bres.) Corner"     # H25 Auktor
$ws.Cells.Item(25, 16).Value = "Storåsens sydsluttning  söder om myren, Hls"  # P25 Lokalnamn
$ws.Cells.Item(25, 17).Value = 599447               # Q25 Ost
$ws.Cells.Item(25, 18).Value = 6820628              # R25 Nord
$ws.Cells.Item(25, 29).Value = "Där stigen delar sig ned mot myren"  # AC25 Publik kommentar (new)
$ws.Cells.Item(25, 34).Value = "Skogsmark"          # AH25 Biotop (new)
$ws.Cells.Item(25, 35).Value = "Barrblandskog kalkpåverkad."  # AI25 Biotop-beskrivning (new)

# --- Row 27: becomes the old row-28 record ---
$ws.Cells.Item(27, 1).Value  = 112306136            # A27 Id
$ws.Cells.Item(27, 2).Value  = 90832                # B27 Taxonsorteringsordning
$ws.Cells.Item(27, 5).Value  = 4368                 # E27 TaxonId
$ws.Cells.Item(27, 6).Value  = "Dofttaggsvamp"      # F27 Artnamn
$ws.Cells.Item(27, 7).Value  = "Hydnellum suaveolens"   # G27 Vetenskapligt namn
$ws.Cells.Item(27, 8).Value  = "(Scop.:Fr.) P. Karst."  # H27 Auktor
$ws.Cells.Item(27, 16).Value = "Bässe söder om stigen, Hls"  # P27 Lokalnamn
$ws.Cells.Item(27, 17).Value = 599416               # Q27 Ost
$ws.Cells.Item(27, 18).Value = 6820643              # R27 Nord
$ws.Cells.Item(27, 29).Value = "Förekommer på flera platser"  # AC27 Publik kommentar

# --- Row 28: becomes the old row-25 record, with its extra fields cleared ---
$ws.Cells.Item(28, 1).Value  = 112306119            # A28 Id
$ws.Cells.Item(28, 2).Value  = 90826                # B28 Taxonsorteringsordning
$ws.Cells.Item(28, 4).Value  = "LC"                 # D28 Rodlistade
$ws.Cells.Item(28, 5).Value  = 4366                 # E28 TaxonId
$ws.Cells.Item(28, 6).Value  = "Skarp dropptaggsvamp"   # F28 Artnamn
$ws.Cells.Item(28, 7).Value  = "Hydnellum peckii"   # G28 Vetenskapligt namn
$ws.Cells.Item(28, 8).Value  = "Banker"             # H28 Auktor
$ws.Cells.Item(28, 29).ClearContents()              # AC28 Publik kommentar (removed)
$ws.Cells.Item(28, 34).ClearContents()              # AH28 Biotop (removed)
$ws.Cells.Item(28, 35).ClearContents()              # AI28 Biotop-beskrivning (removed)
